$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.186.24'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.843.54'
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.56'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6876'
$ws.Range("E6").Value = '  -2.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.0000'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3000'
$ws.Range("E8").Value = '  -2.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07471'
$ws.Range("E9").Value = '  -3.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.24'
$ws.Range("E10").Value = '  -1.54%  '
$ws.Range("E11").Value = '  -1.93%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.842.41'
$ws.Range("E12").Value = '  -0.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.065'
$ws.Range("E13").Value = '  -1.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6838'
$ws.Range("E14").Value = '  -0.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '87.33'
$ws.Range("E15").Value = '  -6.29%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.175'
$ws.Range("E16").Value = '  -6.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.182.31'
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008179'
$ws.Range("E18").Value = '  -1.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.085.19'
$ws.Range("E19").Value = '  -0.72%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '229.45'
$ws.Range("E20").Value = '  -5.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.55'
$ws.Range("E21").Value = '  -1.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9995'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.402'
$ws.Range("E23").Value = '  -1.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9998'
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1450'
$ws.Range("E25").Value = '  -3.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.47'
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.775'
$ws.Range("E27").Value = '  -0.81%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.11'
$ws.Range("E28").Value = '  -1.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.513'
$ws.Range("E29").Value = '  -1.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.282'
$ws.Range("E30").Value = '  +1.31%  '
$ws.Range("E31").Value = '  -0.67%  '
$ws.Range("E32").Value = '  +0.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05260'
$ws.Range("E33").Value = '  +2.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7601'
$ws.Range("E34").Value = '  -3.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.856'
$ws.Range("E35").Value = '  -2.15%  '
$ws.Range("E36").Value = '  -1.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.684'
$ws.Range("E37").Value = '  -0.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.306.18'
$ws.Range("E38").Value = '  -0.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01835'
$ws.Range("E39").Value = '  -1.68%  '
$ws.Range("E40").Value = '  +0.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9346'
$ws.Range("E41").Value = '  -2.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.974'
$ws.Range("E42").Value = '  -1.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '104.98'
$ws.Range("E43").Value = '  -1.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9994'
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.986.34'
$ws.Range("E45").Value = '  -0.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.98'
$ws.Range("E46").Value = '  +0.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5195'
$ws.Range("E47").Value = '  +0.26%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000122'
$ws.Range("E48").Value = '  -0.94%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.528'
$ws.Range("E49").Value = '  -1.79%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.773'
$ws.Range("E50").Value = '  +0.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05958'
$ws.Range("E51").Value = '  +0.82%  '
